# Clean up the header row (row 1) labels: strip stray leading
# tab/space whitespace that had crept into the shared strings and fix
# the mis-encoded Greek letters (gamma/alpha) used as column headers.
# Column order stays the same (A=LC, B=p0x, C=p0y, D=gamma, E=alpha,
# F=v0, G=vspinz, H=vspiny, I=Description); only the text content of
# B1:I1 changes (A1 "LC" is already clean and is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "LC"
$ws.Range("B1").Value = "p0x"
$ws.Range("C1").Value = "p0y"
$ws.Range("D1").Value = "gamma"
$ws.Range("E1").Value = "alpha"
$ws.Range("F1").Value = "v0"
$ws.Range("G1").Value = "vspinz"
$ws.Range("H1").Value = "vspiny"
$ws.Range("I1").Value = "Description"

# Reset the sheet view: scroll back to the top-left (clear the
# previously scrolled topLeftCell="A18") and move the selection from
# I44 to A2.
$ws.Range("A2").Select() | Out-Null
